$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Zeitplanung")
$ws2 = $wb.Worksheets.Item("Ist Arbeitszeit - Übersicht")

# ---------------------------------------------------------------------
# Raw input-cell edits on "Zeitplanung" (all downstream SUM() formulas
# recalc automatically). These cells carry style protection="0" (unlocked)
# so they are editable even though the sheet is protected.
# ---------------------------------------------------------------------
$ws1.Range("Q21").ClearContents()
$ws1.Range("W23").ClearContents()
$ws1.Range("W24").ClearContents()
$ws1.Range("W26").Value = 1
$ws1.Range("X26").ClearContents()
$ws1.Range("AR34").Value = 4
$ws1.Range("AS34").Value = 4
$ws1.Range("C37").Value = 0
$ws1.Range("AR40").Value = 4
$ws1.Range("AS40").Value = 2
$ws1.Range("AS41").Value = 2

# ---------------------------------------------------------------------
# View / selection state: make "Ist Arbeitszeit - Übersicht" the active
# (selected) tab, update each sheet's remembered selection.
# ---------------------------------------------------------------------
[void]$ws1.Range("AY40").Select()
[void]$ws2.Select()
[void]$ws2.Range("P38").Select()

# ---------------------------------------------------------------------
# Move / resize the embedded chart on "Ist Arbeitszeit - Übersicht".
# Values below were solved so the emitted twoCellAnchor lands exactly on
# col4/154078,row1/201707 -> col12/403412,row19/134471 (EMU).
# ---------------------------------------------------------------------
$co = $ws2.ChartObjects(1)
$co.Left   = 281.100836614174
$co.Top    = 31.632401575802316
$co.Width  = 487.132559056118
$co.Height = 276.70586614173175

# ---------------------------------------------------------------------
# Page setup changes on "Ist Arbeitszeit - Übersicht".
# ---------------------------------------------------------------------
$ws2.PageSetup.Orientation = 2
$ws2.PageSetup.Zoom = 95
